# Logged Week 15 and simulated Week 16
# Add a new player, "J.Horsted", as an additional column (U) on both the
# "Rushing" and "Receiving" sheets, matching the header style used by the
# other player columns and using the same placeholder value ("n") used for
# the other players in the data row.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Copy the formatting (font, border, alignment) of the last existing
    # header cell (T1) onto the new header cell (U1) so the new column
    # matches the other headers.
    $ws.Range("T1").Copy()
    $ws.Range("U1").PasteSpecial(-4122)

    $ws.Range("U1").Value = "J.Horsted"
    $ws.Range("U2").Value = "n"
}
